$d = $word.ActiveDocument

# Swap the {{paga}} and {{paga_notula}} placeholders:
#   "Competenze concordate" row currently shows {{paga}}        -> should show {{paga_notula}}
#   "Totale dovuto"         row currently shows {{paga_notula}} -> should show {{paga}}
# A temporary marker avoids the second replace clobbering the first.
$d.Content.Find.Execute("{{paga_notula}}", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{{__SWAP_TMP__}}", 2) | Out-Null

$d.Content.Find.Execute("{{paga}}", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{{paga_notula}}", 2) | Out-Null

$d.Content.Find.Execute("{{__SWAP_TMP__}}", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{{paga}}", 2) | Out-Null
